# EIA Table 1.2.C - roll the report forward from "October 2016" to "November 2016":
#  - update the title and rolling-12-months caption
#  - insert a new "November" monthly data row in the Year-2016 block
#  - refresh the Annual Totals (2014/2015/2016) and the Rolling-12-Months
#    (2015/2016) blocks with the new totals that include November

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title / caption text -------------------------------------------------
$ws.Range("A1").Value = "Table 1.2.C. Net Generation by Energy Source:  Commercial Sector, 2006-November 2016"

# ---- Insert the new "November" row into the Year 2016 monthly block ------
# (row 53 currently holds the "Year to Date" section header; pushing it, and
#  everything below it, down by one row)
$ws.Rows(53).Insert()

# Copy the formatting of the row above (October, row 52) onto the blank new
# row 53 so it matches the other monthly data rows (style only, no values).
$ws.Range("A52:P52").Copy()
$ws.Range("A53:P53").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(53,1).Value = "November"
$ws.Cells.Item(53,2).Value = 39
$ws.Cells.Item(53,3).Value = "NM"
$ws.Cells.Item(53,4).Value = 0.31
$ws.Cells.Item(53,5).Value = 591
$ws.Cells.Item(53,6).Value = 0
$ws.Cells.Item(53,7).Value = 0
$ws.Cells.Item(53,8).Value = "NM"
$ws.Cells.Item(53,9).Value = 38
$ws.Cells.Item(53,10).Value = 197
$ws.Cells.Item(53,11).Value = 0
$ws.Cells.Item(53,12).Value = 84
$ws.Cells.Item(53,13).Value = 960
$ws.Cells.Item(53,14).Value = 467
$ws.Cells.Item(53,15).Value = 505
$ws.Cells.Item(53,16).Value = 505

# ---- Annual Totals block (now rows 55-57: 2014, 2015, 2016) --------------
$ws.Cells.Item(55,2).Value = 551
$ws.Cells.Item(55,3).Value = 236
$ws.Cells.Item(55,4).Value = 8
$ws.Cells.Item(55,5).Value = 6625
$ws.Cells.Item(55,6).Value = 0
$ws.Cells.Item(55,7).Value = 0
$ws.Cells.Item(55,8).Value = 36
$ws.Cells.Item(55,9).Value = 350
$ws.Cells.Item(55,10).Value = 2632
$ws.Cells.Item(55,11).Value = 0
$ws.Cells.Item(55,12).Value = 1074
$ws.Cells.Item(55,13).Value = 11512
$ws.Cells.Item(55,14).Value = 4817
$ws.Cells.Item(55,15).Value = 5167
$ws.Cells.Item(55,16).Value = 5167

$ws.Cells.Item(56,2).Value = 468
$ws.Cells.Item(56,3).Value = 176
$ws.Cells.Item(56,4).Value = 8
$ws.Cells.Item(56,5).Value = 6853
$ws.Cells.Item(56,6).Value = 0
$ws.Cells.Item(56,7).Value = 0
$ws.Cells.Item(56,8).Value = 31
$ws.Cells.Item(56,9).Value = 393
$ws.Cells.Item(56,10).Value = 2562
$ws.Cells.Item(56,11).Value = 0
$ws.Cells.Item(56,12).Value = 1072
$ws.Cells.Item(56,13).Value = 11562
$ws.Cells.Item(56,14).Value = 5340
$ws.Cells.Item(56,15).Value = 5733
$ws.Cells.Item(56,16).Value = 5733

$ws.Cells.Item(57,2).Value = 391
$ws.Cells.Item(57,3).Value = 96
$ws.Cells.Item(57,4).Value = 4
$ws.Cells.Item(57,5).Value = 7148
$ws.Cells.Item(57,6).Value = 0
$ws.Cells.Item(57,7).Value = 0
$ws.Cells.Item(57,8).Value = 54
$ws.Cells.Item(57,9).Value = 534
$ws.Cells.Item(57,10).Value = 2387
$ws.Cells.Item(57,11).Value = 0
$ws.Cells.Item(57,12).Value = 993
$ws.Cells.Item(57,13).Value = 11608
$ws.Cells.Item(57,14).Value = 6710
$ws.Cells.Item(57,15).Value = 7244
$ws.Cells.Item(57,16).Value = 7244

# ---- "Rolling 12 Months Ending in ..." block (now rows 58-60) ------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Cells.Item(59,2).Value = 512
$ws.Cells.Item(59,3).Value = 187
$ws.Cells.Item(59,4).Value = 9
$ws.Cells.Item(59,5).Value = 7455
$ws.Cells.Item(59,6).Value = 0
$ws.Cells.Item(59,7).Value = 0
$ws.Cells.Item(59,8).Value = 33
$ws.Cells.Item(59,9).Value = 413
$ws.Cells.Item(59,10).Value = 2791
$ws.Cells.Item(59,11).Value = 0
$ws.Cells.Item(59,12).Value = 1169
$ws.Cells.Item(59,13).Value = 12569
$ws.Cells.Item(59,14).Value = 5669
$ws.Cells.Item(59,15).Value = 6082
$ws.Cells.Item(59,16).Value = 6082

$ws.Cells.Item(60,2).Value = 432
$ws.Cells.Item(60,3).Value = "NM"
$ws.Cells.Item(60,4).Value = 5
$ws.Cells.Item(60,5).Value = 7765
$ws.Cells.Item(60,6).Value = 0
$ws.Cells.Item(60,7).Value = 0
$ws.Cells.Item(60,8).Value = "NM"
$ws.Cells.Item(60,9).Value = 558
$ws.Cells.Item(60,10).Value = 2629
$ws.Cells.Item(60,11).Value = 0
$ws.Cells.Item(60,12).Value = 1091
$ws.Cells.Item(60,13).Value = 12641
$ws.Cells.Item(60,14).Value = 7059
$ws.Cells.Item(60,15).Value = 7617
$ws.Cells.Item(60,16).Value = 7617

Write-Output "Applied November 2016 update"
